# Regenerate merged AHB files
# 1) Rename header labels in row 1: "_old" suffix -> "_FV2210", "_new" suffix -> "_FV2304"
# 2) Turn the data range A1:U89 into an Excel Table ("Table1")
# 3) Freeze the header row (split pane under row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Create the Excel Table over the used range, with header row
$dataRange = $ws.Range("A1:U89")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze panes beneath the header row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$wb.Save()
